$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated symbol list values for the GitHub Actions price refresh (2023-01-10
# 06:10:48 UTC): each cryptocurrency row's Price (D), Volume 1h % change (E)
# and Hora (G) columns get refreshed figures. Every row's Hora bumps 5 -> 6;
# Price/Volume only change for rows where the scraped market data moved.
#
# {row -> {column letter -> new text value}}
$updates = @{
    2 = @{ "D" = "274.33"; "E" = "-1.47%"; "G" = "6" }
    3 = @{ "D" = "26.64"; "E" = "-2.01%"; "G" = "6" }
    4 = @{ "D" = "4.857"; "E" = "-0.11%"; "G" = "6" }
    5 = @{ "D" = "0.06324"; "E" = "1.18%"; "G" = "6" }
    6 = @{ "D" = "6.891"; "E" = "0.47%"; "G" = "6" }
    7 = @{ "D" = "3.324"; "E" = "1.89%"; "G" = "6" }
    8 = @{ "D" = "1.286"; "E" = "36.57%"; "G" = "6" }
    9 = @{ "D" = "0.8694"; "E" = "-1.15%"; "G" = "6" }
    10 = @{ "D" = "0.1460"; "E" = "0.78%"; "G" = "6" }
    11 = @{ "D" = "0.04990"; "E" = "-3.79%"; "G" = "6" }
    12 = @{ "D" = "0.07375"; "E" = "0.59%"; "G" = "6" }
    13 = @{ "D" = "0.02945"; "E" = "-6.92%"; "G" = "6" }
    14 = @{ "D" = "0.09032"; "E" = "-0.07%"; "G" = "6" }
    15 = @{ "D" = "0.001572"; "E" = "1.25%"; "G" = "6" }
    16 = @{ "D" = "0.0006321"; "E" = "0.83%"; "G" = "6" }
    17 = @{ "D" = "0.005967"; "E" = "0.99%"; "G" = "6" }
    18 = @{ "D" = "3.449"; "E" = "-0.10%"; "G" = "6" }
    19 = @{ "D" = "2.295"; "E" = "1.29%"; "G" = "6" }
    20 = @{ "D" = "0.3127"; "E" = "1.28%"; "G" = "6" }
    21 = @{ "E" = "0.96%"; "G" = "6" }
    22 = @{ "D" = "3.896"; "E" = "1.12%"; "G" = "6" }
    23 = @{ "D" = "0.04360"; "E" = "0.81%"; "G" = "6" }
    24 = @{ "E" = "0.05%"; "G" = "6" }
    25 = @{ "D" = "0.004255"; "E" = "-0.42%"; "G" = "6" }
    26 = @{ "E" = "-0.05%"; "G" = "6" }
    27 = @{ "E" = "0.14%"; "G" = "6" }
    28 = @{ "G" = "6" }
    29 = @{ "G" = "6" }
    30 = @{ "G" = "6" }
    31 = @{ "G" = "6" }
    32 = @{ "G" = "6" }
    33 = @{ "G" = "6" }
    34 = @{ "G" = "6" }
    35 = @{ "G" = "6" }
    36 = @{ "G" = "6" }
    37 = @{ "G" = "6" }
    38 = @{ "G" = "6" }
    39 = @{ "G" = "6" }
    40 = @{ "D" = "0.04047"; "E" = "0.57%"; "G" = "6" }
    41 = @{ "D" = "0.006680"; "E" = "-0.40%"; "G" = "6" }
    42 = @{ "D" = "0.1166"; "E" = "1.19%"; "G" = "6" }
    43 = @{ "D" = "0.002108"; "E" = "0.07%"; "G" = "6" }
    44 = @{ "D" = "0.01222"; "E" = "-8.48%"; "G" = "6" }
    45 = @{ "D" = "0.00005286"; "E" = "4.47%"; "G" = "6" }
    46 = @{ "E" = "-38.75%"; "G" = "6" }
    47 = @{ "G" = "6" }
    48 = @{ "G" = "6" }
    49 = @{ "G" = "6" }
    50 = @{ "G" = "6" }
    51 = @{ "G" = "6" }
}

# All of D/E/G are stored as plain text in this sheet (e.g. "274.33",
# "-1.47%", "6" are literal strings, not numeric cells), so force each cell
# to Text format before writing the new value -- otherwise Excel's
# auto-detection would coerce "274.33" into a number and "-1.47%" into a
# percentage value, changing the cell type. Resetting the style back to
# "Normal" afterwards keeps formatting identical to the original (unstyled)
# cells.
foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $cell = $ws.Range("$col$row")
        $cell.NumberFormat = "@"
        $cell.Value = $cols[$col]
        $cell.Style = "Normal"
    }
}
